$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 31-40: the footnote-marker column (currently F) and its explanation
# column (currently G) both need to shift one column to the left, landing
# in E and F respectively. Use Cut (per source column) so values, shared
# string typing, and cell formatting all move together intact.
$ws.Range("F31:F40").Cut($ws.Range("E31:E40"))
$ws.Range("G31:G40").Cut($ws.Range("F31:F40"))

# Rows 37 and 39 never had an explanation (column G) to begin with, so the
# second Cut above leaves a blank placeholder behind in F37/F39 - drop it
# so those cells go back to being genuinely empty.
$ws.Range("F37").ClearContents()
$ws.Range("F39").ClearContents()

# Update the saved selection to match the new active cell.
$ws.Range("Q24").Select()
